# Refresh the cryptos list (prices + 1h volume %) per the Sun Jul 14 2024
# GitHub Actions data pull. Also RenderToken/NEARProtocol swapped rank
# positions (rows 31/32) between runs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several "Price" cells are digit strings that look like plain numbers
# (e.g. "1.92", "0.999") but must stay literal text, matching the rest of
# the column (prices use "." as a thousands separator, e.g. "59.735.32").
# A direct Range.Value = "1.92" gets auto-parsed by Excel into the number
# 1.92 (losing the trailing zero / the text type). To store it as text
# without flipping the cell to a "Text" number format (which would add a
# style the original file does not have), build the string with a
# text-result formula in a scratch cell, then paste-special VALUES ONLY
# into the destination -- that lands as plain text with the default style.
function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

$ws.Range("D2").Value = '59.735.32'
$ws.Range("E2").Value = '  +3.35%  '
$ws.Range("D3").Value = '3.187.91'
$ws.Range("E3").Value = '  +2.32%  '
$ws.Range("E4").Value = '  -0.04%  '
Set-TextValue $ws.Range("D5") '536.90'
$ws.Range("E5").Value = '  +0.75%  '
Set-TextValue $ws.Range("D6") '144.73'
$ws.Range("E6").Value = '  +4.60%  '
$ws.Range("E7").Value = '  -0.14%  '
Set-TextValue $ws.Range("D8") '0.518'
$ws.Range("E8").Value = '  +4.31%  '
$ws.Range("E9").Value = '  -0.88%  '
$ws.Range("E10").Value = '  +5.08%  '
Set-TextValue $ws.Range("D11") '0.428'
$ws.Range("E11").Value = '  +3.81%  '
$ws.Range("D12").Value = '3.736.30'
$ws.Range("E12").Value = '  +2.21%  '
$ws.Range("E13").Value = '  -0.32%  '
Set-TextValue $ws.Range("D14") '25.95'
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("E15").Value = '  +4.32%  '
$ws.Range("D16").Value = '59.731.59'
$ws.Range("E16").Value = '  +3.15%  '
$ws.Range("D17").Value = '3.220.07'
$ws.Range("E17").Value = '  +3.23%  '
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("E19").Value = '  +1.79%  '
Set-TextValue $ws.Range("D20") '8.22'
$ws.Range("E20").Value = '  +1.46%  '
Set-TextValue $ws.Range("D21") '380.20'
$ws.Range("E21").Value = '  +1.81%  '
$ws.Range("E22").Value = '  -0.06%  '
Set-TextValue $ws.Range("D23") '0.529'
$ws.Range("E23").Value = '  +4.22%  '
Set-TextValue $ws.Range("D24") '70.13'
$ws.Range("E24").Value = '  +1.09%  '
$ws.Range("E25").Value = '  +2.76%  '
$ws.Range("E26").Value = '  +16.27%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").Value = '0.0₃0904'
$ws.Range("E28").Value = '  +2.36%  '
Set-TextValue $ws.Range("D29") '1.92'
$ws.Range("E29").Value = '  +2.54%  '
Set-TextValue $ws.Range("D30") '22.36'
$ws.Range("E30").Value = '  +4.00%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D31") '5.41'
$ws.Range("E31").Value = '  +5.29%  '
$ws.Range("B32").Value = 'RenderToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D32") '6.14'
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("E33").Value = '  +1.97%  '
Set-TextValue $ws.Range("D34") '6.46'
$ws.Range("E34").Value = '  +4.55%  '
Set-TextValue $ws.Range("D35") '156.70'
$ws.Range("E35").Value = '  -2.31%  '
$ws.Range("E36").Value = '  +3.33%  '
Set-TextValue $ws.Range("D37") '25.65'
$ws.Range("E37").Value = '  +0.42%  '
$ws.Range("D38").Value = '2.755.51'
$ws.Range("E38").Value = '  +7.63%  '
Set-TextValue $ws.Range("D39") '0.0712'
$ws.Range("E39").Value = '  +6.22%  '
Set-TextValue $ws.Range("D40") '1.69'
$ws.Range("E40").Value = '  +2.84%  '
Set-TextValue $ws.Range("D41") '4.28'
$ws.Range("E41").Value = '  +3.00%  '
$ws.Range("E42").Value = '  +4.17%  '
Set-TextValue $ws.Range("D43") '39.46'
$ws.Range("E43").Value = '  +2.78%  '
Set-TextValue $ws.Range("D44") '0.0289'
$ws.Range("E44").Value = '  +7.10%  '
$ws.Range("D45").Value = '3.233.56'
$ws.Range("E45").Value = '  +2.41%  '
$ws.Range("E46").Value = '  +2.18%  '
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("E48").Value = '  +5.25%  '
Set-TextValue $ws.Range("D49") '20.52'
$ws.Range("E49").Value = '  +3.04%  '
Set-TextValue $ws.Range("D50") '0.777'
$ws.Range("E50").Value = '  +3.97%  '
Set-TextValue $ws.Range("D51") '0.999'
$ws.Range("E51").Value = '  -0.06%  '
